$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: logistic_embeddings
$ws.Range("C5").Value = 0.37
$ws.Range("D5").Value = 0.488
$ws.Range("E5").Value = 0.523
$ws.Range("F5").Value = 0.5629999999999999
$ws.Range("G5").Value = 0.532
$ws.Range("H5").Value = 0.554

# Row 7: classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.37
$ws.Range("D7").Value = 0.488
$ws.Range("E7").Value = 0.523

# Row 8: BERT-base
$ws.Range("C8").Value = 0.374
$ws.Range("D8").Value = 0.5639999999999999
$ws.Range("E8").Value = 0.6
$ws.Range("F8").Value = 0.637
$ws.Range("G8").Value = 0.621
$ws.Range("H8").Value = 0.65

# Row 9: BERT-base-nli
$ws.Range("B9").Value = 0.384
$ws.Range("C9").Value = 0.487
$ws.Range("E9").Value = 0.623
$ws.Range("G9").Value = 0.597
$ws.Range("H9").Value = 0.626
